# Updates recomputed TPM-based NATMI metrics for the Efnb2-Ephb6 LR pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.71598933333333
$ws.Range("H2").Value = 137.147968
$ws.Range("I2").Value = 0.6549002937372808
$ws.Range("J2").Value = 0.6549002937372808
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08324533333333334
$ws.Range("N2").Value = 0.249736
$ws.Range("O2").Value = 0.05078606388889115
$ws.Range("P2").Value = 0.05078606388889115
$ws.Range("Q2").Value = 3.805642770716445
$ws.Range("R2").Value = 34.250784936448
$ws.Range("S2").Value = 0.03325980815859512
$ws.Range("T2").Value = 0.03325980815859512

# Row 3
$ws.Range("G3").Value = 45.71598933333333
$ws.Range("H3").Value = 137.147968
$ws.Range("I3").Value = 0.6549002937372808
$ws.Range("J3").Value = 0.6549002937372808
$ws.Range("O3").Value = 0.2796082573516313
$ws.Range("P3").Value = 0.2796082573516313
$ws.Range("Q3").Value = 20.95238460595911
$ws.Range("R3").Value = 188.571461453632
$ws.Range("S3").Value = 0.1831155298709525
$ws.Range("T3").Value = 0.1831155298709525

# Row 4
$ws.Range("G4").Value = 45.71598933333333
$ws.Range("H4").Value = 137.147968
$ws.Range("I4").Value = 0.6549002937372808
$ws.Range("J4").Value = 0.6549002937372808
$ws.Range("M4").Value = 1.097575666666667
$ws.Range("N4").Value = 3.292727
$ws.Range("O4").Value = 0.6696056787594775
$ws.Range("P4").Value = 0.6696056787594775
$ws.Range("Q4").Value = 50.17675746985956
$ws.Range("R4").Value = 451.590817228736
$ws.Range("S4").Value = 0.4385249557077331
$ws.Range("T4").Value = 0.4385249557077331

# Row 5
$ws.Range("I5").Value = 0.1818108415648851
$ws.Range("J5").Value = 0.1818108415648851
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08324533333333334
$ws.Range("N5").Value = 0.249736
$ws.Range("O5").Value = 0.05078606388889115
$ws.Range("P5").Value = 0.05078606388889115
$ws.Range("Q5").Value = 1.056507565282667
$ws.Range("R5").Value = 9.508568087543999
$ws.Range("S5").Value = 0.00923345701540732
$ws.Range("T5").Value = 0.009233457015407318

# Row 6
$ws.Range("I6").Value = 0.1818108415648851
$ws.Range("J6").Value = 0.1818108415648851
$ws.Range("O6").Value = 0.2796082573516313
$ws.Range("P6").Value = 0.2796082573516313
$ws.Range("S6").Value = 0.05083581257759105
$ws.Range("T6").Value = 0.05083581257759104

# Row 7
$ws.Range("I7").Value = 0.1818108415648851
$ws.Range("J7").Value = 0.1818108415648851
$ws.Range("M7").Value = 1.097575666666667
$ws.Range("N7").Value = 3.292727
$ws.Range("O7").Value = 0.6696056787594775
$ws.Range("P7").Value = 0.6696056787594775
$ws.Range("Q7").Value = 13.92987389047033
$ws.Range("R7").Value = 125.368865014233
$ws.Range("S7").Value = 0.1217415719718867
$ws.Range("T7").Value = 0.1217415719718867

# Row 8
$ws.Range("G8").Value = 11.24784666666667
$ws.Range("H8").Value = 33.74354
$ws.Range("I8").Value = 0.161130015850732
$ws.Range("J8").Value = 0.161130015850732
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08324533333333334
$ws.Range("N8").Value = 0.249736
$ws.Range("O8").Value = 0.05078606388889115
$ws.Range("P8").Value = 0.05078606388889115
$ws.Range("Q8").Value = 0.9363307450488889
$ws.Range("R8").Value = 8.42697670544
$ws.Range("S8").Value = 0.008183159279413318
$ws.Range("T8").Value = 0.008183159279413318

# Row 9
$ws.Range("G9").Value = 11.24784666666667
$ws.Range("H9").Value = 33.74354
$ws.Range("I9").Value = 0.161130015850732
$ws.Range("J9").Value = 0.161130015850732
$ws.Range("O9").Value = 0.2796082573516313
$ws.Range("P9").Value = 0.2796082573516313
$ws.Range("Q9").Value = 5.155071842162222
$ws.Range("R9").Value = 46.39564657945999
$ws.Range("S9").Value = 0.0450532829390639
$ws.Range("T9").Value = 0.04505328293906389

# Row 10
$ws.Range("G10").Value = 11.24784666666667
$ws.Range("H10").Value = 33.74354
$ws.Range("I10").Value = 0.161130015850732
$ws.Range("J10").Value = 0.161130015850732
$ws.Range("M10").Value = 1.097575666666667
$ws.Range("N10").Value = 3.292727
$ws.Range("O10").Value = 0.6696056787594775
$ws.Range("P10").Value = 0.6696056787594775
$ws.Range("Q10").Value = 12.34536280373111
$ws.Range("R10").Value = 111.10826523358
$ws.Range("S10").Value = 0.1078935736322548
$ws.Range("T10").Value = 0.1078935736322548

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1507006666666667
$ws.Range("H11").Value = 0.452102
$ws.Range("I11").Value = 0.00215884884710222
$ws.Range("J11").Value = 0.00215884884710222
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.08324533333333334
$ws.Range("N11").Value = 0.249736
$ws.Range("O11").Value = 0.05078606388889115
$ws.Range("P11").Value = 0.05078606388889115
$ws.Range("Q11").Value = 0.01254512723022222
$ws.Range("R11").Value = 0.112906145072
$ws.Range("S11").Value = 0.0001096394354753923
$ws.Range("T11").Value = 0.0001096394354753923

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1507006666666667
$ws.Range("H12").Value = 0.452102
$ws.Range("I12").Value = 0.00215884884710222
$ws.Range("J12").Value = 0.00215884884710222
$ws.Range("O12").Value = 0.2796082573516313
$ws.Range("P12").Value = 0.2796082573516313
$ws.Range("Q12").Value = 0.06906857697755556
$ws.Range("R12").Value = 0.621617192798
$ws.Range("S12").Value = 0.0006036319640238301
$ws.Range("T12").Value = 0.0006036319640238298

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1507006666666667
$ws.Range("H13").Value = 0.452102
$ws.Range("I13").Value = 0.00215884884710222
$ws.Range("J13").Value = 0.00215884884710222
$ws.Range("M13").Value = 1.097575666666667
$ws.Range("N13").Value = 3.292727
$ws.Range("O13").Value = 0.6696056787594775
$ws.Range("P13").Value = 0.6696056787594775
$ws.Range("Q13").Value = 0.1654053846837778
$ws.Range("R13").Value = 1.488648462154
$ws.Range("S13").Value = 0.001445577447602998
$ws.Range("T13").Value = 0.001445577447602997
